# Updated symbol list on Mon Jan 23 13:30:20 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto rows that changed, keeping the cells stored as plain text so the
# visual layout/format matches the original scraped sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric/percent-looking strings (e.g. "305.03",
    # "0.76%") are not auto-converted into numbers by Excel, matching the
    # original inline-string cell contents.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    # Reset the style back to the sheet's default so no stray "Text" /
    # quote-prefix formatting is left behind on the cell.
    $cell.Style = "Normal"
}

Set-TextValue "D2" "305.03"
Set-TextValue "E2" "0.76%"
Set-TextValue "D3" "35.91"
Set-TextValue "E3" "-3.21%"
Set-TextValue "D4" "5.069"
Set-TextValue "E4" "1.03%"
Set-TextValue "D5" "0.07947"
Set-TextValue "E5" "1.08%"
Set-TextValue "D6" "2.107"
Set-TextValue "E6" "-4.73%"
Set-TextValue "E7" "3.39%"
Set-TextValue "D8" "7.952"
Set-TextValue "E8" "-0.65%"
Set-TextValue "D9" "0.9222"
Set-TextValue "E9" "0.32%"
Set-TextValue "D10" "0.09675"
Set-TextValue "E10" "1.23%"
Set-TextValue "D11" "0.1845"
Set-TextValue "E11" "-1.88%"
Set-TextValue "D12" "0.08695"
Set-TextValue "E12" "1.45%"
Set-TextValue "D13" "0.03581"
Set-TextValue "E13" "-0.35%"
Set-TextValue "D14" "0.09962"
Set-TextValue "E14" "-0.07%"
Set-TextValue "D15" "0.001436"
Set-TextValue "E15" "-3.12%"
Set-TextValue "D16" "0.005716"
Set-TextValue "E16" "0.51%"
Set-TextValue "D17" "3.464"
Set-TextValue "E17" "0.08%"
Set-TextValue "D18" "2.750"
Set-TextValue "E18" "22.26%"
Set-TextValue "E19" "-0.94%"
Set-TextValue "E20" "2.31%"
Set-TextValue "D21" "5.195"
Set-TextValue "E21" "9.22%"
Set-TextValue "D22" "0.2213"
Set-TextValue "E22" "0.58%"
Set-TextValue "D23" "0.04553"
Set-TextValue "E23" "-0.79%"
Set-TextValue "D24" "0.001239"
Set-TextValue "E24" "0.82%"
Set-TextValue "D25" "0.004884"
Set-TextValue "E25" "9.49%"
Set-TextValue "D26" "0.0001304"
Set-TextValue "E26" "-6.84%"
Set-TextValue "D27" "0.0004761"
Set-TextValue "E27" "0.22%"
Set-TextValue "D39" "0.01845"
Set-TextValue "E39" "2.29%"
Set-TextValue "D40" "0.04737"
Set-TextValue "E40" "0.42%"
Set-TextValue "D41" "0.007927"
Set-TextValue "E41" "-2.45%"
Set-TextValue "D42" "0.1397"
Set-TextValue "E42" "0.11%"
Set-TextValue "D43" "0.007837"
Set-TextValue "E43" "3.75%"
Set-TextValue "D44" "0.002195"
Set-TextValue "E44" "-1.51%"
Set-TextValue "D45" "0.01129"
Set-TextValue "E45" "8.44%"
Set-TextValue "D46" "0.00006296"
Set-TextValue "E46" "2.21%"
Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "0.32%"
Set-TextValue "E48" "0.25%"
Set-TextValue "D49" "49.45"
Set-TextValue "E49" "597.77%"
Set-TextValue "D50" "0.002005"
Set-TextValue "E50" "-25.47%"
Set-TextValue "D51" "0.00002106"
Set-TextValue "E51" "0.32%"
